$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift values in columns B:G down by one row, for rows 2..10 -> 3..11
# (row 11's old data is discarded; row 2 will receive brand new values)
for ($r = 10; $r -ge 2; $r--) {
    $src = $ws.Range("B" + $r + ":G" + $r)
    $dst = $ws.Range("B" + ($r + 1) + ":G" + ($r + 1))
    $dst.Value2 = $src.Value2
}

# New values for row 2 (columns B:G)
$ws.Range("B2").Value = -0.02907897629796788
$ws.Range("C2").Value = 0.3131278957257717
$ws.Range("D2").Value = 0.181524606355785
$ws.Range("E2").Value = 0.4260570458938391
$ws.Range("F2").Value = 0.43998257208981
$ws.Range("G2").Value = 15
